# Applies the "feature selection" re-ranking edit described by the commit:
# both worksheets (final_fail, final_gifted) get several boolean cells
# flipped and every row re-sorted in descending order of the Total (col J)
# count. Cheapest reliable way to land that in one shot is to just restate
# the full post-edit A1:J36 grid for each sheet, cell by cell (bulk
# Range.Value array assignment is not supported by this COM host).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('final_fail')

$rowVals = @('Feature', 'RFE', 'RFECV', 'Logistics', 'Random Forest', 'LightGBM', 'Lasso', 'Ridge', 'Elastic', 'Total')
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks (% of course total)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(2, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Days with no interaction', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(3, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Largest period of inactivity (h)', $true, $true, $false, $true, $true, $true, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(4, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Average session duration (min)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(5, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Submissions (% of course total)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(6, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('On/off campus click ratio', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(7, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Total time online (min)', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(8, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 1 (%)', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(9, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 2 (%)', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(10, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 3 (%)', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(11, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Quizzes started', $true, $true, $true, $true, $false, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(12, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks per session', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(13, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks per day', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(14, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Number of days', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(15, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Resources viewed', $true, $false, $false, $true, $true, $false, $true, $false, 4)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(16, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks on course', $true, $false, $false, $true, $true, $false, $true, $false, 4)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(17, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Number of clicks', $false, $false, $true, $true, $true, $false, $true, $false, 4)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(18, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks on campus', $false, $false, $true, $true, $true, $false, $true, $false, 4)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(19, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Assignments viewed', $true, $true, $false, $false, $true, $false, $true, $false, 4)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(20, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Days with no interaction (%)', $false, $false, $false, $true, $true, $false, $true, $false, 3)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(21, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Files downloaded', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(22, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Forum posts', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(23, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Discussions viewed', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(24, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Links viewed', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(25, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Assignments submitted', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(26, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks on folder', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(27, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 5 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(28, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 4 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(29, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Number of sessions', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(30, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 10 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(31, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 9 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(32, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 8 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(33, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 7 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(34, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 6 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(35, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks on forum', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(36, $c).Value = $rowVals[$c - 1]
}

$ws = $wb.Worksheets.Item('final_gifted')

$rowVals = @('Feature', 'RFE', 'RFECV', 'Logistics', 'Random Forest', 'LightGBM', 'Lasso', 'Ridge', 'Elastic', 'Total')
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('On/off campus click ratio', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(2, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Days with no interaction', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(3, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Total time online (min)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(4, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks per day', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(5, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks per session', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(6, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks (% of course total)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(7, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Number of clicks', $false, $true, $true, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(8, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 4 (%)', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(9, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Resources viewed', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(10, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks on course', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(11, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 3 (%)', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(12, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 1 (%)', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(13, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Average session duration (min)', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(14, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Largest period of inactivity (h)', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(15, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks on campus', $true, $true, $false, $true, $true, $false, $true, $false, 5)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(16, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 2 (%)', $true, $false, $false, $true, $true, $false, $true, $false, 4)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(17, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Number of days', $true, $true, $false, $false, $true, $false, $true, $false, 4)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(18, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Assignments viewed', $true, $true, $false, $true, $false, $false, $true, $false, 4)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(19, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Submissions (% of course total)', $false, $false, $true, $false, $true, $false, $true, $false, 3)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(20, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Days with no interaction (%)', $false, $false, $false, $true, $true, $false, $true, $false, 3)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(21, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks on forum', $false, $true, $true, $false, $false, $false, $true, $false, 3)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(22, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Clicks on folder', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(23, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Quizzes started', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(24, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Forum posts', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(25, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Number of sessions', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(26, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 8 (%)', $false, $true, $false, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(27, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 7 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(28, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 6 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(29, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 5 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(30, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Links viewed', $false, $false, $true, $false, $false, $false, $true, $false, 2)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(31, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Assignments submitted', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(32, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Discussions viewed', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(33, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Files downloaded', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(34, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 10 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(35, $c).Value = $rowVals[$c - 1]
}
$rowVals = @('Start of Session 9 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(36, $c).Value = $rowVals[$c - 1]
}

